$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "date" column (F) values forward by one week (7 days),
# keeping the existing date formatting/style on each cell.
$ws.Range("F2").Value = 44522
$ws.Range("F3").Value = 44521
$ws.Range("F4").Value = 44520
$ws.Range("F5").Value = 44519
$ws.Range("F6").Value = 44518
$ws.Range("F7").Value = 44517
